$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.993.75'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').Value = '1.873.13'
$ws.Range('E3').Value = '  +0.71%  '
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '305.68'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').Value = '0.9995'
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('D7').Value = '0.5063'
$ws.Range('E7').Value = '  +0.55%  '
$ws.Range('D8').Value = '0.3656'
$ws.Range('E8').Value = '  -1.42%  '
$ws.Range('D9').Value = '0.07187'
$ws.Range('E9').Value = '  +0.95%  '
$ws.Range('D10').Value = '0.8951'
$ws.Range('E10').Value = '  +1.49%  '
$ws.Range('D11').Value = '20.72'
$ws.Range('E11').Value = '  +0.99%  '
$ws.Range('D12').Value = '1.864.57'
$ws.Range('E12').Value = '  +0.40%  '
$ws.Range('D13').Value = '0.07516'
$ws.Range('E13').Value = '  -0.62%  '
$ws.Range('D14').Value = '95.34'
$ws.Range('E14').Value = '  +7.10%  '
$ws.Range('D15').Value = '5.244'
$ws.Range('E15').Value = '  -0.83%  '
$ws.Range('D16').Value = '1.0000'
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('E17').Value = '  +1.98%  '
$ws.Range('E18').Value = '  +1.35%  '
$ws.Range('D19').Value = '0.9990'
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('D20').Value = '27.042.13'
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').Value = '2.103.89'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('E23').Value = '  -0.42%  '
$ws.Range('D24').Value = '6.429'
$ws.Range('E24').Value = '  -0.38%  '
$ws.Range('D25').Value = '148.32'
$ws.Range('E25').Value = '  +0.99%  '
$ws.Range('D26').Value = '1.790'
$ws.Range('E26').Value = '  -2.98%  '
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('D28').Value = '2.073'
$ws.Range('E28').Value = '  -1.05%  '
$ws.Range('D29').Value = '113.29'
$ws.Range('E29').Value = '  +0.77%  '
$ws.Range('D30').Value = '4.696'
$ws.Range('E30').Value = '  +1.08%  '
$ws.Range('D31').Value = '4.687'
$ws.Range('E31').Value = '  +0.56%  '
$ws.Range('D32').Value = '0.09165'
$ws.Range('E32').Value = '  +1.58%  '
$ws.Range('D33').Value = '0.05145'
$ws.Range('E33').Value = '  +0.60%  '
$ws.Range('D34').Value = '0.7537'
$ws.Range('E34').Value = '  +4.67%  '
$ws.Range('D35').Value = '2.987'
$ws.Range('E35').Value = '  -1.14%  '
$ws.Range('D36').Value = '1.159'
$ws.Range('E36').Value = '  +1.25%  '
$ws.Range('D37').Value = '3.231'
$ws.Range('E37').Value = '  +6.58%  '
$ws.Range('D38').Value = '2.584'
$ws.Range('E38').Value = '  +5.48%  '
$ws.Range('D39').Value = '0.5655'
$ws.Range('E39').Value = '  +7.39%  '
$ws.Range('D40').Value = '0.02001'
$ws.Range('E40').Value = '  -1.58%  '
$ws.Range('D41').Value = '1.074'
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').Value = '6.599'
$ws.Range('E42').Value = '  +1.46%  '
$ws.Range('D43').Value = '116.04'
$ws.Range('E43').Value = '  +0.90%  '
$ws.Range('D44').Value = '8.565'
$ws.Range('E44').Value = '  +4.23%  '
$ws.Range('D45').Value = '0.1476'
$ws.Range('E45').Value = '  +0.79%  '
$ws.Range('D46').Value = '0.4733'
$ws.Range('E46').Value = '  +2.96%  '
$ws.Range('D47').Value = '0.9992'
$ws.Range('E47').Value = '  -0.05%  '
$ws.Range('D48').Value = '10.09'
$ws.Range('E48').Value = '  +1.65%  '
$ws.Range('E49').Value = '  +0.33%  '
$ws.Range('E50').Value = '  +1.22%  '
$ws.Range('D51').Value = '63.29'
$ws.Range('E51').Value = '  -0.73%  '
